# Template formatting fix: shrink the "LOCAL DO TREINAMENTO" merge-field
# placeholder text (shape "CaixaDeTexto 6" on slide 1) from 9pt to 7pt and
# resize the text box to match the new autofit height.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(3)                 # "CaixaDeTexto 6"

# The shape has two paragraphs:
#   1) "LOCAL DO TREINAMENTO: " (bold, 6pt)  -> unchanged
#   2) "[local_treinamento]" (9pt)           -> shrink to 7pt
$tr    = $shp.TextFrame.TextRange
$para2 = $tr.Paragraphs(2, 1)
$para2.Font.Size = 7

# The text box uses <a:spAutoFit/>, so shrinking the font reduces the
# rendered/autofit height of the box. Match the recalculated height.
$shp.Height = 290934 / 12700
